$wb = $excel.ActiveWorkbook

# --- Sheet1: fill in the "입력" (answer) column D for the questionnaire rows ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2,4).Value  = "심화"
$ws1.Cells.Item(3,4).Value  = "컴퓨터"
$ws1.Cells.Item(4,4).Value  = "신입학생"
$ws1.Cells.Item(5,4).Value  = 2016
$ws1.Cells.Item(6,4).Value  = 2016

$ws1.Cells.Item(9,4).Value  = 21

$ws1.Cells.Item(12,4).Value = 4

$ws1.Cells.Item(16,4).Value = 28

$ws1.Cells.Item(18,4).Value = 42
$ws1.Cells.Item(19,4).Value = 84

$ws1.Cells.Item(25,4).Value = "O"
$ws1.Cells.Item(26,4).Value = "O"
$ws1.Cells.Item(27,4).Value = 2
$ws1.Cells.Item(28,4).Value = 4
$ws1.Cells.Item(29,4).Value = "X"
$ws1.Cells.Item(30,4).Value = "X"
$ws1.Cells.Item(31,4).Value = 140
$ws1.Cells.Item(32,4).Value = 2.5

$ws1.Cells.Item(35,4).Value = "X"
$ws1.Cells.Item(36,4).Value = 2
$ws1.Cells.Item(37,4).Value = "O"
$ws1.Cells.Item(38,4).Value = 3
$ws1.Cells.Item(39,4).Value = "O"

# --- Sheet2 (수학필수): move the cell selection, without disturbing the active tab ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F10").Select()

# --- Add a new trailing, empty worksheet named "Sheet2" ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet2"

# --- Restore Sheet1 as the active tab with its new selection ---
$ws1.Activate()
$ws1.Range("D40").Select()
